$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.6200000000001
$ws.Range("G2").Value = 0.0003897078414272226
$ws.Range("H2").Value = 0.005305557989998039
$ws.Range("K2").Value = 4.47320923398759
$ws.Range("L2").Value = "[1.63300824410737, 7.31341022386781]"
$ws.Range("M2").Value = 0.002092836964294742
$ws.Range("N2").Value = 0.002193985764704864
$ws.Range("O2").Value = -1.308210754648002
$ws.Range("P2").Value = "[-2.0503687789194647, -0.5660527303765388]"
$ws.Range("Q2").Value = 0.000584403642466258
$ws.Range("R2").Value = 0.001168807284932516
$ws.Range("S2").Value = 14.25275997232328
$ws.Range("T2").Value = "[12.702230074550286, 15.803289870096272]"
$ws.Range("W2").Value = 4.709669669669687
$ws.Range("X2").Value = 2.037837837837842
$ws.Range("Y2").Value = 7.381501501501531

# Row 3 updates
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 23.01000000000016
$ws.Range("G3").Value = 0.00091623248289463
$ws.Range("H3").Value = 0.005305557989998039
$ws.Range("K3").Value = 4.468396279404046
$ws.Range("L3").Value = "[1.6207190401994627, 7.316073518608629]"
$ws.Range("M3").Value = 0.002193985764704864
$ws.Range("N3").Value = 0.002193985764704864
$ws.Range("O3").Value = -0.2012631930227693
$ws.Range("P3").Value = "[-0.9937370155499243, 0.5912106295043857]"
$ws.Range("Q3").Value = 0.6176976900069806
$ws.Range("R3").Value = 0.6176976900069806
$ws.Range("S3").Value = 13.33069447056868
$ws.Range("T3").Value = "[11.671776839715168, 14.9896121014222]"
$ws.Range("W3").Value = 0.7370570570570614
$ws.Range("X3").Value = -2.165105105105121
$ws.Range("Y3").Value = 3.639219219219243
